$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.538.50"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.927.05"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'376.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.47%  "
$ws.Range("D6").Value = "'104.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("D10").Value = "'36.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.0838"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "'18.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "3.389.71"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "2.929.46"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "'0.939"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.12%  "
$ws.Range("D18").Value = "51.461.66"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'3.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.08%  "
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").Value = "'68.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'261.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("E26").Value = "  -5.35%  "
$ws.Range("D27").Value = "'4.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.83%  "
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'7.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("D34").Value = "'51.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("D36").Value = "'34.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0426"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("E39").Value = "  -8.73%  "
$ws.Range("D40").Value = "'16.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -7.53%  "
$ws.Range("E42").Value = "  -5.76%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "'124.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'21.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.59%  "
$ws.Range("E46").Value = "  -6.19%  "
$ws.Range("D47").Value = "'0.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.29%  "
$ws.Range("D48").Value = "2.022.30"
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'3.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").Value = "3.217.01"
$ws.Range("E51").Value = "  -2.70%  "
